# Append two new resale-number rows (15 and 16) dated 2023-06-02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ RowNum = 15; A = "2023-06-02"; B = "10:16:34"; C = "Friday"; D = "22";
       E = 120492; F = 133982; G = 158003; H = 129556; I = 174025; J = 111624;
       K = 198997; L = 217488; M = 170999; N = 118375; O = 37878; P = 34893;
       Q = 49993;  R = -1;     S = 36370;  T = -1 },
    @{ RowNum = 16; A = "2023-06-02"; B = "16:32:08"; C = "Friday"; D = "22";
       E = 120621; F = 133845; G = 158538; H = 130038; I = 174274; J = 112011;
       K = 199359; L = 217948; M = 171156; N = 118531; O = 37968; P = 34867;
       Q = 50103;  R = -1;     S = 36674;  T = -1 }
)

foreach ($row in $newRows) {
    $rn = $row.RowNum

    # Columns A-D hold text (date/time/weekday/week are stored as plain text
    # strings, not real Excel dates/numbers). Pre-format the range as Text so
    # the values aren't auto-converted to a date serial / number, then clear
    # the formatting back to the workbook default so no stray number-format
    # style gets attached to the cells.
    $textRange = $ws.Range(("A{0}:D{0}" -f $rn))
    $textRange.NumberFormat = "@"
    $ws.Range("A$rn").Value = $row.A
    $ws.Range("B$rn").Value = $row.B
    $ws.Range("C$rn").Value = $row.C
    $ws.Range("D$rn").Value = $row.D
    $textRange.ClearFormats()

    # Columns E-T are plain numbers.
    $ws.Range("E$rn").Value = $row.E
    $ws.Range("F$rn").Value = $row.F
    $ws.Range("G$rn").Value = $row.G
    $ws.Range("H$rn").Value = $row.H
    $ws.Range("I$rn").Value = $row.I
    $ws.Range("J$rn").Value = $row.J
    $ws.Range("K$rn").Value = $row.K
    $ws.Range("L$rn").Value = $row.L
    $ws.Range("M$rn").Value = $row.M
    $ws.Range("N$rn").Value = $row.N
    $ws.Range("O$rn").Value = $row.O
    $ws.Range("P$rn").Value = $row.P
    $ws.Range("Q$rn").Value = $row.Q
    $ws.Range("R$rn").Value = $row.R
    $ws.Range("S$rn").Value = $row.S
    $ws.Range("T$rn").Value = $row.T
}
